$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Columns B:I keep width 19.140625, column J gets widened + "best fit"
$ws.Range("J1:J16384").ColumnWidth = 123.7109375

# --- New cell values in column J (rows 2-7 existing table, plus new rows 8-14) ---
$ws.Range("J2").Value = "Fix serial errors"
$ws.Range("J3").Value = "SerBEAT = serial(str); %<--change this appropriately"
$ws.Range("J4").Value = "set(SerBEAT,'BaudRate', 9600, 'DataBits', 8, 'Parity', 'none','StopBits', 1, 'FlowControl', 'none');"
$ws.Range("J5").Value = "fopen(SerBEAT); %--open the serial port to the PIC"
$ws.Range("J6").Value = "NOTE: try using 'serial' call with static members. Ex:"
$ws.Range("J7").Value = "s = serial('COM1','baudrate',4800);`nOR: s2 = serial('COM2','BaudRate',1200,'DataBits',7);"
$ws.Range("J8").Value = "query for member values:"
$ws.Range("J9").Value = "get(s1,{'Type','Name','Port'})"
$ws.Range("J10").Value = "ans = "
$ws.Range("J11").Value = "    'serial'    'Serial-COM1'    'COM1'"
$ws.Range("J12").Value = "UPDATE: fixed by using"
$ws.Range("J13").Value = "upper(COM Port String)"
$ws.Range("J14").Value = "to make string capitals. Don't ask why it needs it, just accept it."

# --- Fonts / styles for the new cells ---
# J3, J4: Lucida Console 10pt, left/vcenter
$ws.Range("J3:J4").Font.Name = "Lucida Console"
$ws.Range("J3:J4").Font.Size = 10
$ws.Range("J3:J4").HorizontalAlignment = -4131  # xlLeft
$ws.Range("J3:J4").VerticalAlignment = -4108    # xlCenter

# J5: Lucida Console 10pt, vcenter only
$ws.Range("J5").Font.Name = "Lucida Console"
$ws.Range("J5").Font.Size = 10
$ws.Range("J5").VerticalAlignment = -4108       # xlCenter

# J9:J11: Consolas 10pt, dark gray, left/vcenter, indent 2
$ws.Range("J9:J11").Font.Name = "Consolas"
$ws.Range("J9:J11").Font.Size = 10
$ws.Range("J9:J11").Font.Color = 4210752
$ws.Range("J9:J11").HorizontalAlignment = -4131 # xlLeft
$ws.Range("J9:J11").VerticalAlignment = -4108   # xlCenter
$ws.Range("J9:J11").IndentLevel = 2

# J12: green fill, wrap text
$ws.Range("J12").Interior.Color = 5287936
$ws.Range("J12").WrapText = $true

# J8, J13, J14: wrap text (style 1, same as default wrap style)
$ws.Range("J8").WrapText = $true
$ws.Range("J13").WrapText = $true
$ws.Range("J14").WrapText = $true

# --- Row height adjustments ---
$ws.Rows.Item(7).RowHeight = 47.25

# --- View changes ---
$ws.Range("J16").Select()
$excel.ActiveWindow.ScrollColumn = 7
